$wb = $excel.ActiveWorkbook

# "展览" (Exhibitions) sheet — F column "想去人数" (want-to-go count) updates
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F6").Value = 5340
$wsExh.Range("F8").Value = 5373
$wsExh.Range("F9").Value = 626
$wsExh.Range("F10").Value = 8
$wsExh.Range("F11").Value = 1372
$wsExh.Range("F12").Value = 3

# "全部类型" (All types) sheet — same events, shifted one row down
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 5340
$wsAll.Range("F9").Value = 5373
$wsAll.Range("F10").Value = 626
$wsAll.Range("F11").Value = 8
$wsAll.Range("F12").Value = 1372
$wsAll.Range("F13").Value = 3
